# katalog.xlsx edit: add a "Steher-Abstand" (post spacing) input row to the
# Draht_Matten sheet and parametrize its price formula, and remove the now
# redundant "Montage (€/m)" row from Brix_Zaun_Stab.

$wb = $excel.ActiveWorkbook

# --- Draht_Matten: insert new row 4 "Steher-Abstand" -----------------------
$wsMatten = $wb.Worksheets.Item("Draht_Matten")
$wsMatten.Activate()

$wsMatten.Rows.Item(4).Insert()

$wsMatten.Range("A4").Value = "Auswahl"
$wsMatten.Range("B4").Value = "Steher-Abstand"
$wsMatten.Range("C4").Value = "Dist"
$wsMatten.Range("D4").Value = "Standard (2.5m):2.5, Verkürzt (2.0m):2.0, Eng (1.25m):1.25"

# Update the price formula (row 10 -> now row 11) to use the new Dist
# variable instead of the hard-coded 2.5 spacing.
$wsMatten.Range("E11").Value = "(L * P_Matte * F_Faktor) + ((math.ceil(L/Dist)+1) * ((P_Saeule * F_Faktor) + (Ist_Beton * 2 * P_Sack) + ((1-Ist_Beton) * P_Konsole))) + (L * P_Arbeit)"

$wsMatten.Columns.Item(2).ColumnWidth = 26.28515625
$wsMatten.Columns.Item(3).ColumnWidth = 17.7109375

$wsMatten.Range("E11").Select()

# --- Brix_Zaun_Stab: remove the separate "Montage (€/m)" row ---------------
$wsZaun = $wb.Worksheets.Item("Brix_Zaun_Stab")
$wsZaun.Activate()

$wsZaun.Rows.Item(7).Delete()

$wsZaun.Columns.Item(1).ColumnWidth = 14.85546875
$wsZaun.Columns.Item(2).ColumnWidth = 21.5703125
$wsZaun.Columns.Item(3).ColumnWidth = 12.7109375
$wsZaun.Columns.Item(4).ColumnWidth = 42.5703125
$wsZaun.Columns.Item(5).ColumnWidth = 19.42578125

$wsZaun.Range("D9").Select()

# --- Re-activate Draht_Matten as the last-active sheet ---------------------
$wsMatten.Activate()
